# Auto-generated Excel COM-interop script
# Commit: Update automàtic: dades i banners [2026-02-05 09:17]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-05 09:15:38'
$ws.Range('K2').Value = '0.1 MJ/m2'
$ws.Range('O2').Value = '-2.6 °C'
$ws.Range('E3').Value = '2026-02-05 09:15:40'
$ws.Range('K3').Value = '0.0 MJ/m2'
$ws.Range('N3').Value = '-4.9 °C 7:20 TU'
$ws.Range('E4').Value = '2026-02-05 09:15:43'
$ws.Range('E5').Value = '2026-02-05 09:15:45'
$ws.Range('E6').Value = '2026-02-05 09:15:48'
$ws.Range('E7').Value = '2026-02-05 09:15:50'
$ws.Range('M7').Value = '10.6 °C 7:01 TU'
$ws.Range('O7').Value = '9.7 °C'
$ws.Range('E8').Value = '2026-02-05 09:15:53'
$ws.Range('E9').Value = '2026-02-05 09:15:55'
$ws.Range('M9').Value = '1.3 °C 8:29 TU'
$ws.Range('O9').Value = '-0.2 °C'
$ws.Range('E10').Value = '2026-02-05 09:15:58'
$ws.Range('E11').Value = '2026-02-05 09:16:00'
$ws.Range('E12').Value = '2026-02-05 09:16:02'
$ws.Range('E13').Value = '2026-02-05 09:16:05'
$ws.Range('E14').Value = '2026-02-05 09:16:07'
$ws.Range('E15').Value = '2026-02-05 09:16:10'
$ws.Range('J15').Value = '994.0 hPa'
$ws.Range('K15').Value = '0.2 MJ/m2'
$ws.Range('L15').Value = '9.7 km/h - 162º 8:29 TU'
$ws.Range('M15').Value = '6.6 °C 8:29 TU'
$ws.Range('O15').Value = '1.7 °C'
$ws.Range('E16').Value = '2026-02-05 09:16:12'
$ws.Range('H16').Value = '98%'
$ws.Range('K16').Value = '0.1 MJ/m2'
$ws.Range('L16').Value = '20.5 km/h - 267º 7:44 TU'
$ws.Range('M16').Value = '3.1 °C 8:29 TU'
$ws.Range('O16').Value = '2.2 °C'
$ws.Range('E17').Value = '2026-02-05 09:16:15'
$ws.Range('E18').Value = '2026-02-05 09:16:18'
$ws.Range('E19').Value = '2026-02-05 09:16:20'
$ws.Range('E20').Value = '2026-02-05 09:16:22'
$ws.Range('E21').Value = '2026-02-05 09:16:25'
$ws.Range('J21').Value = '995.1 hPa'
$ws.Range('K21').Value = '0.2 MJ/m2'
$ws.Range('L21').Value = '8.6 km/h - 252º 8:19 TU'
$ws.Range('M21').Value = '4.7 °C 8:25 TU'
$ws.Range('O21').Value = '0.9 °C'
$ws.Range('E22').Value = '2026-02-05 09:16:28'
$ws.Range('I22').Value = '0.1 mm'
$ws.Range('K22').Value = '0.2 MJ/m2'
$ws.Range('M22').Value = '5.4 °C 8:29 TU'
$ws.Range('O22').Value = '3.3 °C'
$ws.Range('E23').Value = '2026-02-05 09:16:30'
$ws.Range('J23').Value = '993.7 hPa'
$ws.Range('K23').Value = '0.2 MJ/m2'
$ws.Range('E24').Value = '2026-02-05 09:16:33'
$ws.Range('E25').Value = '2026-02-05 09:16:35'
$ws.Range('J25').Value = '997.6 hPa'
$ws.Range('K25').Value = '0.3 MJ/m2'
$ws.Range('M25').Value = '0.7 °C 8:25 TU'
$ws.Range('O25').Value = '-0.7 °C'
$ws.Range('E26').Value = '2026-02-05 09:16:38'
$ws.Range('E27').Value = '2026-02-05 09:16:40'
$ws.Range('E28').Value = '2026-02-05 09:16:43'
$ws.Range('J28').Value = '997.8 hPa'
$ws.Range('L28').Value = '20.9 km/h - 262º 8:13 TU'
$ws.Range('M28').Value = '0.6 °C 8:29 TU'
$ws.Range('O28').Value = '-1.2 °C'
$ws.Range('E29').Value = '2026-02-05 09:16:45'
$ws.Range('E30').Value = '2026-02-05 09:16:48'
$ws.Range('E31').Value = '2026-02-05 09:16:50'
$ws.Range('E32').Value = '2026-02-05 09:16:53'
$ws.Range('E33').Value = '2026-02-05 09:16:55'
$ws.Range('E34').Value = '2026-02-05 09:16:58'
$ws.Range('K34').Value = '0.1 MJ/m2'
$ws.Range('O34').Value = '0.5 °C'
$ws.Range('E35').Value = '2026-02-05 09:17:00'
$ws.Range('E36').Value = '2026-02-05 09:17:03'
